# Update gh-pages to output generated at 456a3b4
#
# Source data refresh: the "想去人数" (want-to-go headcount) column (F) on
# the 展览 / 演出 / 全部类型 sheets moved slightly since the last scrape.
# Apply the new counts cell by cell (numbers only, no formatting changes).

$wb = $excel.ActiveWorkbook

# ---- Sheet: 展览 (column F = "想去人数") ----
$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 576
$ws.Cells.Item(3, 6).Value = 5443
$ws.Cells.Item(5, 6).Value = 456
$ws.Cells.Item(7, 6).Value = 1008
$ws.Cells.Item(8, 6).Value = 378
$ws.Cells.Item(9, 6).Value = 1345
$ws.Cells.Item(12, 6).Value = 3068
$ws.Cells.Item(14, 6).Value = 119
$ws.Cells.Item(16, 6).Value = 189
$ws.Cells.Item(17, 6).Value = 22
$ws.Cells.Item(18, 6).Value = 136
$ws.Cells.Item(20, 6).Value = 971
$ws.Cells.Item(21, 6).Value = 349
$ws.Cells.Item(23, 6).Value = 3530
$ws.Cells.Item(24, 6).Value = 1111
$ws.Cells.Item(25, 6).Value = 2799
$ws.Cells.Item(27, 6).Value = 1961
$ws.Cells.Item(28, 6).Value = 4043
$ws.Cells.Item(29, 6).Value = 108
$ws.Cells.Item(30, 6).Value = 914
$ws.Cells.Item(31, 6).Value = 462
$ws.Cells.Item(32, 6).Value = 1286
$ws.Cells.Item(33, 6).Value = 50
$ws.Cells.Item(34, 6).Value = 28
$ws.Cells.Item(35, 6).Value = 995
$ws.Cells.Item(36, 6).Value = 1264
$ws.Cells.Item(37, 6).Value = 59
$ws.Cells.Item(38, 6).Value = 1030
$ws.Cells.Item(39, 6).Value = 667
$ws.Cells.Item(40, 6).Value = 525
$ws.Cells.Item(41, 6).Value = 407
$ws.Cells.Item(42, 6).Value = 17

# ---- Sheet: 演出 (column F = "想去人数") ----
$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(5, 6).Value = 24
$ws.Cells.Item(10, 6).Value = 904

# ---- Sheet: 全部类型 (column F = "想去人数") ----
$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 577
$ws.Cells.Item(3, 6).Value = 577
$ws.Cells.Item(4, 6).Value = 5443
$ws.Cells.Item(9, 6).Value = 378
$ws.Cells.Item(10, 6).Value = 1345
$ws.Cells.Item(11, 6).Value = 3068
$ws.Cells.Item(14, 6).Value = 119
$ws.Cells.Item(17, 6).Value = 189
$ws.Cells.Item(18, 6).Value = 904
$ws.Cells.Item(21, 6).Value = 136
$ws.Cells.Item(22, 6).Value = 971
$ws.Cells.Item(23, 6).Value = 349
$ws.Cells.Item(24, 6).Value = 3530
$ws.Cells.Item(27, 6).Value = 1111
$ws.Cells.Item(28, 6).Value = 2799
$ws.Cells.Item(29, 6).Value = 1961
$ws.Cells.Item(30, 6).Value = 4043
$ws.Cells.Item(32, 6).Value = 108
$ws.Cells.Item(33, 6).Value = 914
$ws.Cells.Item(34, 6).Value = 1286
$ws.Cells.Item(35, 6).Value = 28
$ws.Cells.Item(36, 6).Value = 995
$ws.Cells.Item(38, 6).Value = 1264
$ws.Cells.Item(39, 6).Value = 59
$ws.Cells.Item(40, 6).Value = 1030
$ws.Cells.Item(42, 6).Value = 667
$ws.Cells.Item(44, 6).Value = 407
